# Rename a handful of "English Variable Name" values (column C) in the
# codebook worksheet, as described in the commit: "changed some variable
# names". The `<` / `>` characters used to denote "under"/"over" age
# thresholds are spelled out instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "males_under12"
$ws.Range("C12").Value = "males_over12"
$ws.Range("C14").Value = "fem_under12"
$ws.Range("C15").Value = "fem_over12"
$ws.Range("C17").Value = "pers_under12"
$ws.Range("C18").Value = "pers_over12"
$ws.Range("C78").Value = "under10_yrs"

# The author also scrolled the sheet down and moved the selection before
# saving. Reproduce the final cursor/viewport position.
$ws.Range("F84").Select()
$excel.ActiveWindow.ScrollRow = 99
$excel.ActiveWindow.ScrollColumn = 1
